$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ETH-USD"
$ws.Range("B2").Value = "Custom"
$ws.Range("F2").Value = "B"

$ws.Range("F2").Select()
